$d = $word.ActiveDocument

# The document currently ends with a trailing empty paragraph; insert the
# new content (a sub-bullet under "Inferring Return Types" plus a blank
# paragraph) immediately before it, using a full OOXML snippet so the
# formatting (list numbering, fonts, shading) lands exactly as authored.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)

$insertionPoint = $lastPara.Range.Duplicate
$insertionPoint.Collapse(1)

$xml = '<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="5"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Segoe UI"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t>It does this by looking at the types of the values after a function’s </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Courier New"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:shd w:val="clear" w:color="auto" w:fill="EAE9ED"/>
    </w:rPr>
    <w:t>return</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Segoe UI"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
    </w:rPr>
    <w:t> statements.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
    </w:rPr>
  </w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>'
[void]$insertionPoint.InsertXML($xml)

# The original trailing empty paragraph now picks up a left indent.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalPara.LeftIndent = 54
